$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.8096947552348297
$ws.Cells.Item(2, 3).Value = 0.1689056610607338
$ws.Cells.Item(2, 4).Value = 0.3803482496744834
$ws.Cells.Item(2, 6).Value = 0.9510534095220038
$ws.Cells.Item(2, 7).Value = 0.3826376112167651
$ws.Cells.Item(2, 8).Value = 0.5329323714551606
$ws.Cells.Item(2, 10).Value = 0.3464697890742627
$ws.Cells.Item(2, 14).Value = 0.9354102844883698
$ws.Cells.Item(2, 15).Value = 1.773290268889582

$ws.Cells.Item(3, 2).Value = 0.7167389832155493
$ws.Cells.Item(3, 3).Value = 0.1483254934948945
$ws.Cells.Item(3, 4).Value = 0.3705370198418336
$ws.Cells.Item(3, 6).Value = 0.9450809235930251
$ws.Cells.Item(3, 7).Value = 0.3797178100942844
$ws.Cells.Item(3, 8).Value = 0.5355345402302447
$ws.Cells.Item(3, 10).Value = 0.3347488705399257
$ws.Cells.Item(3, 14).Value = 0.9350643962688565
$ws.Cells.Item(3, 15).Value = 1.772263363456148

$ws.Cells.Item(4, 2).Value = 0.6595868745604321
$ws.Cells.Item(4, 3).Value = 0.1356287670299992
$ws.Cells.Item(4, 4).Value = 0.364675798423832
$ws.Cells.Item(4, 6).Value = 0.942048742712096
$ws.Cells.Item(4, 7).Value = 0.3782623373390948
$ws.Cells.Item(4, 8).Value = 0.5374245285424308
$ws.Cells.Item(4, 10).Value = 0.3277567447083527
$ws.Cells.Item(4, 14).Value = 0.9352330977286698
$ws.Cells.Item(4, 15).Value = 1.772949590126231

$ws.Cells.Item(5, 2).Value = 0.636279138145909
$ws.Cells.Item(5, 3).Value = 0.1304399357110242
$ws.Cells.Item(5, 4).Value = 0.3623284075038384
$ws.Cells.Item(5, 6).Value = 0.940972652058683
$ws.Cells.Item(5, 7).Value = 0.3777538240050191
$ws.Cells.Item(5, 8).Value = 0.5382681896068959
$ws.Cells.Item(5, 10).Value = 0.3249587956166522
$ws.Cells.Item(5, 14).Value = 0.9353979254587728
$ws.Cells.Item(5, 15).Value = 1.773559857939617

$ws.Cells.Item(6, 2).Value = 0.632407877490067
$ws.Cells.Item(6, 3).Value = 0.1295774508314196
$ws.Cells.Item(6, 4).Value = 0.3619411117175702
$ws.Cells.Item(6, 6).Value = 0.9408036009450669
$ws.Cells.Item(6, 7).Value = 0.3776744893776183
$ws.Cells.Item(6, 8).Value = 0.5384127162741734
$ws.Cells.Item(6, 10).Value = 0.3244973014978143
$ws.Cells.Item(6, 14).Value = 0.9354311056145121
$ws.Cells.Item(6, 15).Value = 1.773681147240623

$ws.Cells.Item(7, 2).Value = 0.6592726082486706
$ws.Cells.Item(7, 3).Value = 0.1355588481029599
$ws.Cells.Item(7, 4).Value = 0.3646439740755341
$ws.Cells.Item(7, 6).Value = 0.9420335843294865
$ws.Cells.Item(7, 7).Value = 0.3782551370760174
$ws.Cells.Item(7, 8).Value = 0.5374356089682593
$ws.Cells.Item(7, 10).Value = 0.3277188025286222
$ws.Cells.Item(7, 14).Value = 0.9352349312853576
$ws.Cells.Item(7, 15).Value = 1.772956482367832

$ws.Cells.Item(8, 2).Value = 0.7776604616412897
$ws.Cells.Item(8, 3).Value = 0.1618223817689
$ws.Cells.Item(8, 4).Value = 0.3769316042418041
$ws.Cells.Item(8, 6).Value = 0.948862238354053
$ws.Cells.Item(8, 7).Value = 0.3815607280225208
$ws.Cells.Item(8, 8).Value = 0.5337689354549582
$ws.Cells.Item(8, 10).Value = 0.3423859299081613
$ws.Cells.Item(8, 14).Value = 0.9352120787045664
$ws.Cells.Item(8, 15).Value = 1.772662603981075

$ws.Cells.Item(9, 2).Value = 1.009151831127951
$ws.Cells.Item(9, 3).Value = 0.2128313067397869
$ws.Cells.Item(9, 4).Value = 0.4023157481129545
$ws.Cells.Item(9, 6).Value = 0.9672986458918729
$ws.Cells.Item(9, 7).Value = 0.3907298397471806
$ws.Cells.Item(9, 8).Value = 0.528898387985933
$ws.Cells.Item(9, 10).Value = 0.3727750064809214
$ws.Cells.Item(9, 14).Value = 0.9381809390412741
$ws.Cells.Item(9, 15).Value = 1.782558814709915

$ws.Cells.Item(10, 2).Value = 1.178762144603127
$ws.Cells.Item(10, 3).Value = 0.2499906959071723
$ws.Cells.Item(10, 4).Value = 0.4217470403735888
$ws.Cells.Item(10, 6).Value = 0.9839330980578751
$ws.Cells.Item(10, 7).Value = 0.3991203841794828
$ws.Cells.Item(10, 8).Value = 0.5267360843268847
$ws.Cells.Item(10, 10).Value = 0.3961016000720292
$ws.Cells.Item(10, 14).Value = 0.9421876907292841
$ws.Cells.Item(10, 15).Value = 1.79625221844293

$ws.Cells.Item(11, 2).Value = 1.255809325023847
$ws.Cells.Item(11, 3).Value = 0.26682351804169
$ws.Cells.Item(11, 4).Value = 0.43075590098303
$ws.Cells.Item(11, 6).Value = 0.9921744037128093
$ws.Cells.Item(11, 7).Value = 0.4033001243330148
$ws.Cells.Item(11, 8).Value = 0.526060337093142
$ws.Cells.Item(11, 10).Value = 0.4069324212280634
$ws.Cells.Item(11, 14).Value = 0.9444045929241156
$ws.Cells.Item(11, 15).Value = 1.803884808649229

$ws.Cells.Item(12, 2).Value = 1.284968006041481
$ws.Cells.Item(12, 3).Value = 0.2731870941619832
$ws.Cells.Item(12, 4).Value = 0.4341915849327904
$ws.Cells.Item(12, 6).Value = 0.9953923053529365
$ws.Cells.Item(12, 7).Value = 0.4049353118698917
$ws.Cells.Item(12, 8).Value = 0.525848753014543
$ws.Cells.Item(12, 10).Value = 0.4110654150026249
$ws.Cells.Item(12, 14).Value = 0.9453005567784203
$ws.Cells.Item(12, 15).Value = 1.806977481448001

$ws.Cells.Item(13, 2).Value = 1.278688963605987
$ws.Cells.Item(13, 3).Value = 0.2718170643520068
$ws.Cells.Item(13, 4).Value = 0.4334505738721361
$ws.Cells.Item(13, 6).Value = 0.9946949518623427
$ws.Cells.Item(13, 7).Value = 0.404580810274183
$ws.Cells.Item(13, 8).Value = 0.5258923503717909
$ws.Cells.Item(13, 10).Value = 0.4101738949028402
$ws.Cells.Item(13, 14).Value = 0.9451050867430837
$ws.Cells.Item(13, 15).Value = 1.80630240922099

$ws.Cells.Item(14, 2).Value = 1.258208585885541
$ws.Cells.Item(14, 3).Value = 0.2673472690550227
$ws.Cells.Item(14, 4).Value = 0.4310380721190938
$ws.Cells.Item(14, 6).Value = 0.9924371955878684
$ws.Cells.Item(14, 7).Value = 0.4034336004683325
$ws.Cells.Item(14, 8).Value = 0.526042041788827
$ws.Cells.Item(14, 10).Value = 0.4072718116264156
$ws.Cells.Item(14, 14).Value = 0.9444771738616566
$ws.Cells.Item(14, 15).Value = 1.804135185517339

$ws.Cells.Item(15, 2).Value = 1.245661445490498
$ws.Cells.Item(15, 3).Value = 0.2646079913793926
$ws.Cells.Item(15, 4).Value = 0.4295634954915499
$ws.Cells.Item(15, 6).Value = 0.9910669040653062
$ws.Cells.Item(15, 7).Value = 0.4027377336963553
$ws.Cells.Item(15, 8).Value = 0.5261395030903628
$ws.Cells.Item(15, 10).Value = 0.405498317905483
$ws.Cells.Item(15, 14).Value = 0.9440999069821032
$ws.Cells.Item(15, 15).Value = 1.802834070701635

$ws.Cells.Item(16, 2).Value = 1.173724584826175
$ws.Cells.Item(16, 3).Value = 0.2488891638365942
$ws.Cells.Item(16, 4).Value = 0.4211616878580173
$ws.Cells.Item(16, 6).Value = 0.9834080888155796
$ws.Cells.Item(16, 7).Value = 0.398854549688167
$ws.Cells.Item(16, 8).Value = 0.5267864443549684
$ws.Cells.Item(16, 10).Value = 0.3953981996842941
$ws.Cells.Item(16, 14).Value = 0.9420507249011507
$ws.Cells.Item(16, 15).Value = 1.795781699390233

$ws.Cells.Item(17, 2).Value = 1.129564432188943
$ws.Cells.Item(17, 3).Value = 0.2392276488082814
$ws.Cells.Item(17, 4).Value = 0.4160507525305945
$ws.Cells.Item(17, 6).Value = 0.9788824404136705
$ws.Cells.Item(17, 7).Value = 0.3965654414382698
$ws.Cells.Item(17, 8).Value = 0.5272622038701371
$ws.Cells.Item(17, 10).Value = 0.3892583375188394
$ws.Cells.Item(17, 14).Value = 0.9408944295076367
$ws.Cells.Item(17, 15).Value = 1.791815147217733

$ws.Cells.Item(18, 2).Value = 1.104154483252444
$ws.Cells.Item(18, 3).Value = 0.233663926684784
$ws.Cells.Item(18, 4).Value = 0.4131270375720817
$ws.Cells.Item(18, 6).Value = 0.9763428582415941
$ws.Cells.Item(18, 7).Value = 0.3952829451267803
$ws.Cells.Item(18, 8).Value = 0.5275648268430757
$ws.Cells.Item(18, 10).Value = 0.3857475031531266
$ws.Cells.Item(18, 14).Value = 0.9402664765104305
$ws.Cells.Item(18, 15).Value = 1.789665751297747

$ws.Cells.Item(19, 2).Value = 1.095549411104571
$ws.Cells.Item(19, 3).Value = 0.2317790123055374
$ws.Cells.Item(19, 4).Value = 0.4121398634724187
$ws.Cells.Item(19, 6).Value = 0.9754938926277674
$ws.Cells.Item(19, 7).Value = 0.3948545700078654
$ws.Cells.Item(19, 8).Value = 0.5276722659345268
$ws.Cells.Item(19, 10).Value = 0.3845623406577658
$ws.Cells.Item(19, 14).Value = 0.9400602440588699
$ws.Cells.Item(19, 15).Value = 1.788960666777058

$ws.Cells.Item(20, 2).Value = 1.134266423386407
$ws.Cells.Item(20, 3).Value = 0.2402568273377312
$ws.Cells.Item(20, 4).Value = 0.4165931695885376
$ws.Cells.Item(20, 6).Value = 0.9793576353886237
$ws.Cells.Item(20, 7).Value = 0.3968055860368906
$ws.Cells.Item(20, 8).Value = 0.5272085590417248
$ws.Cells.Item(20, 10).Value = 0.3899097986213462
$ws.Cells.Item(20, 14).Value = 0.9410136791817933
$ws.Cells.Item(20, 15).Value = 1.792223721610952

$ws.Cells.Item(21, 2).Value = 1.264224655272017
$ws.Cells.Item(21, 3).Value = 0.2686604485670614
$ws.Cells.Item(21, 4).Value = 0.4317460260929522
$ws.Cells.Item(21, 6).Value = 0.9930977169523203
$ws.Cells.Item(21, 7).Value = 0.4037691395226517
$ws.Cells.Item(21, 8).Value = 0.5259968710661127
$ws.Cells.Item(21, 10).Value = 0.4081233659708232
$ws.Cells.Item(21, 14).Value = 0.9446600761794741
$ws.Cells.Item(21, 15).Value = 1.804766254680203

$ws.Cells.Item(22, 2).Value = 1.349057592136603
$ws.Cells.Item(22, 3).Value = 0.2871615992612249
$ws.Cells.Item(22, 4).Value = 0.4417904219116053
$ws.Cells.Item(22, 6).Value = 1.002643683805218
$ws.Cells.Item(22, 7).Value = 0.4086258216666039
$ws.Cells.Item(22, 8).Value = 0.5254632312205558
$ws.Cells.Item(22, 10).Value = 0.4202111748538613
$ws.Cells.Item(22, 14).Value = 0.9473722679208691
$ws.Cells.Item(22, 15).Value = 1.814143393131133

$ws.Cells.Item(23, 2).Value = 1.303790570609522
$ws.Cells.Item(23, 3).Value = 0.2772930200308679
$ws.Cells.Item(23, 4).Value = 0.4364166736793607
$ws.Cells.Item(23, 6).Value = 0.9974969815312846
$ws.Cells.Item(23, 7).Value = 0.4060056795967597
$ws.Cells.Item(23, 8).Value = 0.5257244027945092
$ws.Cells.Item(23, 10).Value = 0.4137428166546897
$ws.Cells.Item(23, 14).Value = 0.9458946810046029
$ws.Cells.Item(23, 15).Value = 1.809030492122105

$ws.Cells.Item(24, 2).Value = 1.132140720595601
$ws.Cells.Item(24, 3).Value = 0.2397915643156807
$ws.Cells.Item(24, 4).Value = 0.4163478972560597
$ws.Cells.Item(24, 6).Value = 0.9791426057732906
$ws.Cells.Item(24, 7).Value = 0.3966969121987205
$ws.Cells.Item(24, 8).Value = 0.5272327212171604
$ws.Cells.Item(24, 10).Value = 0.3896152137117355
$ws.Cells.Item(24, 14).Value = 0.94095965170483
$ws.Cells.Item(24, 15).Value = 1.792038597024572

$ws.Cells.Item(25, 2).Value = 0.9466048778025424
$ws.Cells.Item(25, 3).Value = 0.1990865465273259
$ws.Cells.Item(25, 4).Value = 0.3953111447165725
$ws.Cells.Item(25, 6).Value = 0.9617695726370101
$ws.Cells.Item(25, 7).Value = 0.3879600455585717
$ws.Cells.Item(25, 8).Value = 0.5299674567558412
$ws.Cells.Item(25, 10).Value = 0.3643789630252314
$ws.Cells.Item(25, 14).Value = 0.9370563794358873
$ws.Cells.Item(25, 15).Value = 1.778756502801286
